# Rebuilds a paragraph's text content from a list of strings, each becoming
# its own <w:r> run (all sharing the paragraph's default run formatting,
# i.e. <w:rPr><w:lang w:val="en-US"/></w:rPr>), mirroring the way the
# target OOXML splits a paragraph's sentence into many same-styled runs.
function Rebuild-Paragraph($d, $paraIndex, $runTexts) {
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    # Keep the paragraph mark, wipe every run (and any proofErr markers)
    # that currently live inside the paragraph.
    $r.MoveEnd(1, -1)
    $r.Delete()

    $fullText = ""
    foreach ($t in $runTexts) {
        $fullText = $fullText + $t
    }
    $r.InsertAfter($fullText)

    # Re-apply the language mark that the rest of the document uses, on the
    # freshly (re)created run that now spans the whole paragraph.
    $p2 = $d.Paragraphs($paraIndex)
    $langRange = $p2.Range
    $langRange.LanguageID = "en-US"

    # Split the single run into one run per entry of $runTexts by toggling
    # Bold on/off across each segment in turn -- Word (and this host) only
    # coalesces adjacent runs when nothing ever forced a boundary between
    # them, so a harmless on/off flip leaves the final formatting untouched
    # while still breaking the run apart at that point.
    $pos = $p2.Range.Start
    foreach ($t in $runTexts) {
        $len = $t.Length
        if ($len -gt 0) {
            $seg = $d.Range($pos, $pos + $len)
            $seg.Bold = 1
            $seg.Bold = 0
        }
        $pos = $pos + $len
    }
}

$d = $word.ActiveDocument
$rsquo = [char]0x2019

# Paragraph: "My 10 years of experience ... Innovation Manager capacity."
# -> "My Ph.D. in Materials Science ... Senior Engineer capacity."
$p1Runs = @(
    "My Ph.D. in ",
    "Materials Science and Engineering",
    " and more than 2 years as a postdoctoral fellow in computational solid and fluid mechanics and advanced fracture mechanics, along with my ",
    "multidisciplinary engineering background",
    ", provide me with the skills to ",
    "drive",
    " innovat",
    ("ion at Amgen" + $rsquo + "s Device Engineering group "),
    "in a ",
    "Senior",
    " Engineer",
    " capacity."
)

# Paragraph: "I offer proficiency in managing international R&D projects..."
# -> "During my professional journey I contributed to projects ... as outlined in my CV."
$p2Runs = @(
    "During my professional journey I contributed to projects in several fields, from ankle biomechanics to the design of nano-sized electromagnetic sensors",
    ",",
    " from modeling and prediction of damage in fiber-reinforced composites to multi-scale modeling of wood, from coupled fluid-structure interaction to large displacement analysis of cracking in hydrogels. ",
    "I offer proficiency in several ",
    "CAD and CAE tools, as well as ",
    "computational methods of solid mechanics (FEM, BEM), fluid mechanics (LBM, FVM), fracture and damage mechanics (CZM, VCCT, J-integral, interaction integrals), mesh generation and computational geometry (Delaunay triangulation, transfinite interpolation, elliptic/parabolic/hyperbolic mesh smoothing",
    "). I have multiple years of experience in mechanical testing and I am expert in several ",
    "programming languages",
    ", ",
    "as outlined in my CV",
    "."
)

# Paragraph: "I am eager to put my skills and experience in service of GrowHub's mission..."
# -> "... in service of Amgen's mission..."
$p3Runs = @(
    "I am eager to put my skills and experience in service of ",
    "Amgen",
    "'s mission, and I am confident we should arrange a time to meet. In the meantime, I wish to thank you for taking the time to consider my application and review my qualifications."
)

Rebuild-Paragraph $d 16 $p1Runs
Rebuild-Paragraph $d 18 $p2Runs
Rebuild-Paragraph $d 22 $p3Runs
